$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.67714804971649
$ws.Range("D2").Value = 3.564411215264732
$ws.Range("E2").Value = 21.99485615319196
$ws.Range("F2").Value = 26.03181424364242
$ws.Range("G2").Value = 34.32052207884342
$ws.Range("H2").Value = 13.97939580066629
$ws.Range("I2").Value = 23.98181680477459
$ws.Range("L2").Value = 9.626826924212519
$ws.Range("M2").Value = 14.85035265206367
$ws.Range("N2").Value = 18.63355544373744
$ws.Range("B3").Value = 14.4205579541273
$ws.Range("D3").Value = 3.585436484294257
$ws.Range("E3").Value = 21.50245760165358
$ws.Range("F3").Value = 25.45776901868672
$ws.Range("G3").Value = 33.10597213900947
$ws.Range("H3").Value = 13.88530575374864
$ws.Range("I3").Value = 24.14141469297891
$ws.Range("L3").Value = 9.517564374321589
$ws.Range("M3").Value = 14.73288718484196
$ws.Range("N3").Value = 18.66667591949177
$ws.Range("B4").Value = 14.26293407330604
$ws.Range("D4").Value = 3.599136223565167
$ws.Range("E4").Value = 21.19193451841554
$ws.Range("F4").Value = 25.10914029813316
$ws.Range("G4").Value = 32.35184039402914
$ws.Range("H4").Value = 13.83176848321417
$ws.Range("I4").Value = 24.24422785080581
$ws.Range("L4").Value = 9.451725828301702
$ws.Range("M4").Value = 14.66264586337133
$ws.Range("N4").Value = 18.68915693536854
$ws.Range("B5").Value = 14.19876283830352
$ws.Range("D5").Value = 3.604917202213941
$ws.Range("E5").Value = 21.06344111101905
$ws.Range("F5").Value = 24.96826872432398
$ws.Range("G5").Value = 32.04299286122564
$ws.Range("H5").Value = 13.81103386781345
$ws.Range("I5").Value = 24.28734038328838
$ws.Range("L5").Value = 9.425237976121394
$ws.Range("M5").Value = 14.63451998932385
$ws.Range("N5").Value = 18.69885858003892
$ws.Range("B6").Value = 14.18811337404117
$ws.Range("D6").Value = 3.605889088541452
$ws.Range("E6").Value = 21.04199016302699
$ws.Range("F6").Value = 24.94495625618017
$ws.Range("G6").Value = 31.99163411798916
$ws.Range("H6").Value = 13.80765673527283
$ws.Range("I6").Value = 24.29457269151644
$ws.Range("L6").Value = 9.420861136994764
$ws.Range("M6").Value = 14.62988049022835
$ws.Range("N6").Value = 18.70050220434124
$ws.Range("B7").Value = 14.26206827964193
$ws.Range("D7").Value = 3.599213385884497
$ws.Range("E7").Value = 21.19020937284005
$ws.Range("F7").Value = 25.10723529804362
$ws.Range("G7").Value = 32.34768057895551
$ws.Range("H7").Value = 13.83148444553928
$ws.Range("I7").Value = 24.24480435610602
$ws.Range("L7").Value = 9.451367183377206
$ws.Range("M7").Value = 14.66226449972178
$ws.Range("N7").Value = 18.68928558546616
$ws.Range("B8").Value = 14.58873988393457
$ws.Range("D8").Value = 3.571496448864896
$ws.Range("E8").Value = 21.82683850977338
$ws.Range("F8").Value = 25.83322579322299
$ws.Range("G8").Value = 33.90382034464381
$ws.Range("H8").Value = 13.94608583437634
$ws.Range("I8").Value = 24.03584797557832
$ws.Range("L8").Value = 9.588910806094683
$ws.Range("M8").Value = 14.80947472792392
$ws.Range("N8").Value = 18.64453086573225
$ws.Range("B9").Value = 15.22524992806816
$ws.Range("D9").Value = 3.523434668454282
$ws.Range("E9").Value = 23.00594521577864
$ws.Range("F9").Value = 27.27728792100064
$ws.Range("G9").Value = 36.86570156089288
$ws.Range("H9").Value = 14.20354939510641
$ws.Range("I9").Value = 23.66416548309616
$ws.Range("L9").Value = 9.867297648084158
$ws.Range("M9").Value = 15.11197979878792
$ws.Range("N9").Value = 18.57373803279007
$ws.Range("B10").Value = 15.68587768091406
$ws.Range("D10").Value = 3.49199006615908
$ws.Range("E10").Value = 23.82460970414975
$ws.Range("F10").Value = 28.33763819945642
$ws.Range("G10").Value = 38.95945525753248
$ws.Range("H10").Value = 14.41136823211733
$ws.Range("I10").Value = 23.41407936594046
$ws.Range("L10").Value = 10.07541837877805
$ws.Range("M10").Value = 15.34112515466744
$ws.Range("N10").Value = 18.53200838153947
$ws.Range("B11").Value = 15.89300093449207
$ws.Range("D11").Value = 3.478532114400577
$ws.Range("E11").Value = 24.18569444438903
$ws.Range("F11").Value = 28.81729983696272
$ws.Range("G11").Value = 39.88916478503279
$ws.Range("H11").Value = 14.50965211245478
$ws.Range("I11").Value = 23.30525256332878
$ws.Range("L11").Value = 10.17051655236447
$ws.Range("M11").Value = 15.44653092193184
$ws.Range("N11").Value = 18.51524329869608
$ws.Range("B12").Value = 15.97101418605263
$ws.Range("D12").Value = 3.473558353039623
$ws.Range("E12").Value = 24.32072297198254
$ws.Range("F12").Value = 28.99834659091529
$ws.Range("G12").Value = 40.23760098312565
$ws.Range("H12").Value = 14.54738019868714
$ws.Range("I12").Value = 23.26474930735136
$ws.Range("L12").Value = 10.20655863722761
$ws.Range("M12").Value = 15.48658448096738
$ws.Range("N12").Value = 18.50921257804604
$ws.Range("B13").Value = 15.95423241673411
$ws.Range("D13").Value = 3.474624081822396
$ws.Range("E13").Value = 24.29171921872858
$ws.Range("F13").Value = 28.9593843030359
$ws.Range("G13").Value = 40.16272540306161
$ws.Range("H13").Value = 14.53923255922707
$ws.Range("I13").Value = 23.27344100867624
$ws.Range("L13").Value = 10.19879545250162
$ws.Range("M13").Value = 15.47795252520961
$ws.Range("N13").Value = 18.51049728250796
$ws.Range("B14").Value = 15.89942794533792
$ws.Range("D14").Value = 3.478120460967675
$ws.Range("E14").Value = 24.19683794174504
$ws.Range("F14").Value = 28.83220761892281
$ws.Range("G14").Value = 39.91790535912504
$ws.Range("H14").Value = 14.51274596890372
$ws.Range("I14").Value = 23.30190618239126
$ws.Range("L14").Value = 10.17348131792216
$ws.Range("M14").Value = 15.44982353942548
$ws.Range("N14").Value = 18.51474078371464
$ws.Range("B15").Value = 15.86580187401022
$ws.Range("D15").Value = 3.480278067271298
$ws.Range("E15").Value = 24.13849600925559
$ws.Range("F15").Value = 28.75422549098423
$ws.Range("G15").Value = 39.76746420808126
$ws.Range("H15").Value = 14.49658773272004
$ws.Range("I15").Value = 23.31943391682811
$ws.Range("L15").Value = 10.15797875463306
$ws.Range("M15").Value = 15.43261092128092
$ws.Range("N15").Value = 18.51738141202592
$ws.Range("B16").Value = 15.67228678046053
$ws.Range("D16").Value = 3.492886670450413
$ws.Range("E16").Value = 23.80077723694764
$ws.Range("F16").Value = 28.30621929137024
$ws.Range("G16").Value = 38.89820761595607
$ws.Range("H16").Value = 14.40501811385308
$ws.Range("I16").Value = 23.42129056741025
$ws.Range("L16").Value = 10.069209396384
$ws.Range("M16").Value = 15.33425762943078
$ws.Range("N16").Value = 18.53314855740237
$ws.Range("B17").Value = 15.55289968469666
$ws.Range("D17").Value = 3.500838944749183
$ws.Range("E17").Value = 23.5906389120503
$ws.Range("F17").Value = 28.03054303251019
$ws.Range("G17").Value = 38.35885535689354
$ws.Range("H17").Value = 14.34978261916836
$ws.Range("I17").Value = 23.48503902350671
$ws.Range("L17").Value = 10.01483887974727
$ws.Range("M17").Value = 15.27419917273433
$ws.Range("N17").Value = 18.54338849157248
$ws.Range("B18").Value = 15.48400782066714
$ws.Range("D18").Value = 3.505492497013754
$ws.Range("E18").Value = 23.46871069509778
$ws.Range("F18").Value = 27.87173887043233
$ws.Range("G18").Value = 38.04651126087594
$ws.Range("H18").Value = 14.31836700915457
$ws.Range("I18").Value = 23.52217047800137
$ws.Range("L18").Value = 9.983608006326877
$ws.Range("M18").Value = 15.23976667009438
$ws.Range("N18").Value = 18.54948706019533
$ws.Range("B19").Value = 15.46064607546937
$ws.Range("D19").Value = 3.507081760994882
$ws.Range("E19").Value = 23.42724791691172
$ws.Range("F19").Value = 27.81793524641418
$ws.Range("G19").Value = 37.94040375517175
$ws.Range("H19").Value = 14.30779194078218
$ws.Range("I19").Value = 23.53482251387069
$ws.Range("L19").Value = 9.973041847738797
$ws.Range("M19").Value = 15.22812845601953
$ws.Range("N19").Value = 18.5515878301867
$ws.Range("B20").Value = 15.56563231616364
$ws.Range("D20").Value = 3.49998416636276
$ws.Range("E20").Value = 23.6131190095865
$ws.Range("F20").Value = 28.05991584407972
$ws.Range("G20").Value = 38.41649273426334
$ws.Range("H20").Value = 14.35562605763025
$ws.Range("I20").Value = 23.4782047852472
$ws.Range("L20").Value = 10.02062262894537
$ws.Range("M20").Value = 15.28058116492902
$ws.Range("N20").Value = 18.5422768278793
$ws.Range("B21").Value = 15.91553729716244
$ws.Range("D21").Value = 3.477090159440617
$ws.Range("E21").Value = 24.22475376870956
$ws.Range("F21").Value = 28.86958008071444
$ws.Range("G21").Value = 39.98991592080361
$ws.Range("H21").Value = 14.5205121112724
$ws.Range("I21").Value = 23.29352610491092
$ws.Range("L21").Value = 10.18091610551638
$ws.Range("M21").Value = 15.4580821722425
$ws.Range("N21").Value = 18.51348574684657
$ws.Range("B22").Value = 16.14174253265511
$ws.Range("D22").Value = 3.462841855912771
$ws.Range("E22").Value = 24.61451468189198
$ws.Range("F22").Value = 29.39521417732294
$ws.Range("G22").Value = 40.99697677455399
$ws.Range("H22").Value = 14.63123403285822
$ws.Range("I22").Value = 23.17694798551674
$ws.Range("L22").Value = 10.28584104612989
$ws.Range("M22").Value = 15.57488422003049
$ws.Range("N22").Value = 18.49652128366523
$ws.Range("B23").Value = 16.02126211892589
$ws.Range("D23").Value = 3.470380824233782
$ws.Range("E23").Value = 24.40742879726007
$ws.Range("F23").Value = 29.11506000500511
$ws.Range("G23").Value = 40.46154101495614
$ws.Range("H23").Value = 14.57187858879444
$ws.Range("I23").Value = 23.23879194138168
$ws.Range("L23").Value = 10.22983553543513
$ws.Range("M23").Value = 15.51248163557305
$ws.Range("N23").Value = 18.50540642543623
$ws.Range("B24").Value = 15.55987668109623
$ws.Range("D24").Value = 3.500370357301628
$ws.Range("E24").Value = 23.60295923044242
$ws.Range("F24").Value = 28.04663735913042
$ws.Range("G24").Value = 38.39044189837084
$ws.Range("H24").Value = 14.3529831778883
$ws.Range("I24").Value = 23.48129304399888
$ws.Range("L24").Value = 10.01800770809114
$ws.Range("M24").Value = 15.27769556533635
$ws.Range("N24").Value = 18.54277875219618
$ws.Range("B25").Value = 15.0539717861268
$ws.Range("D25").Value = 3.535760074752475
$ws.Range("E25").Value = 22.69496595931887
$ws.Range("F25").Value = 26.88582385447188
$ws.Range("G25").Value = 36.07707473395003
$ws.Range("H25").Value = 14.13052118578258
$ws.Range("I25").Value = 23.76066215102051
$ws.Range("L25").Value = 9.791233184519633
$ws.Range("M25").Value = 15.0288203080906
$ws.Range("N25").Value = 18.59107908638953
